# Updates cryptocurrency price/volume data (and a handful of re-ordered rows)
# in the "cryptos" worksheet, per the Oct 28 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the A1 cell reference plus its new text value. Values are
# written as literal text (prices/percentages are display strings in this
# sheet, e.g. "69.733.50" / "  +2.87%  "), so every cell is forced to the
# "@" (text) number format before the assignment -- otherwise Excel would
# auto-coerce number-looking strings (like "596.19") into real numbers --
# and formats are cleared afterward so the cell style matches the original
# (un-styled) cells.
$updates = @(
    @{ Cell = "D2"; Value = "69.733.50" },
    @{ Cell = "E2"; Value = "  +2.87%  " },
    @{ Cell = "D3"; Value = "2.507.03" },
    @{ Cell = "E3"; Value = "  +0.55%  " },
    @{ Cell = "E4"; Value = "  +0.10%  " },
    @{ Cell = "D5"; Value = "596.19" },
    @{ Cell = "E5"; Value = "  +1.54%  " },
    @{ Cell = "D6"; Value = "176.94" },
    @{ Cell = "E6"; Value = "  +0.16%  " },
    @{ Cell = "E7"; Value = "  +0.07%  " },
    @{ Cell = "E8"; Value = "  +0.60%  " },
    @{ Cell = "D9"; Value = "2.502.92" },
    @{ Cell = "E9"; Value = "  +0.38%  " },
    @{ Cell = "D10"; Value = "0.158" },
    @{ Cell = "E10"; Value = "  +11.61%  " },
    @{ Cell = "E11"; Value = "  -0.48%  " },
    @{ Cell = "E12"; Value = "  +0.75%  " },
    @{ Cell = "E13"; Value = "  +1.40%  " },
    @{ Cell = "B14"; Value = "WrappedliquidstakedEther2.0" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth" },
    @{ Cell = "D14"; Value = "2.971.58" },
    @{ Cell = "E14"; Value = "  +0.76%  " },
    @{ Cell = "B15"; Value = "Avalanche" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax" },
    @{ Cell = "D15"; Value = "25.90" },
    @{ Cell = "E15"; Value = "  +0.86%  " },
    @{ Cell = "B16"; Value = "WrappedBTC" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc" },
    @{ Cell = "D16"; Value = "69.519.72" },
    @{ Cell = "E16"; Value = "  +2.69%  " },
    @{ Cell = "B17"; Value = "ShibaInu" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib" },
    @{ Cell = "D17"; Value = "0.0000177" },
    @{ Cell = "E17"; Value = "  +3.50%  " },
    @{ Cell = "D18"; Value = "2.525.74" },
    @{ Cell = "E18"; Value = "  +0.45%  " },
    @{ Cell = "B19"; Value = "BitcoinCash" },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch" },
    @{ Cell = "D19"; Value = "364.15" },
    @{ Cell = "E19"; Value = "  +3.89%  " },
    @{ Cell = "B20"; Value = "Chainlink" },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" },
    @{ Cell = "D20"; Value = "11.02" },
    @{ Cell = "E20"; Value = "  +0.57%  " },
    @{ Cell = "B21"; Value = "Uniswap" },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" },
    @{ Cell = "D21"; Value = "7.53" },
    @{ Cell = "E21"; Value = "  -0.05%  " },
    @{ Cell = "E22"; Value = "  +0.19%  " },
    @{ Cell = "E23"; Value = "  -0.10%  " },
    @{ Cell = "D24"; Value = "70.65" },
    @{ Cell = "E24"; Value = "  -0.18%  " },
    @{ Cell = "E25"; Value = "  -1.15%  " },
    @{ Cell = "D26"; Value = "9.04" },
    @{ Cell = "E26"; Value = "  -0.80%  " },
    @{ Cell = "D27"; Value = "1.67" },
    @{ Cell = "E27"; Value = "  -4.24%  " },
    @{ Cell = "D28"; Value = "2.643.47" },
    @{ Cell = "E28"; Value = "  +1.08%  " },
    @{ Cell = "D29"; Value = "0.959" },
    @{ Cell = "E29"; Value = "  -2.40%  " },
    @{ Cell = "D30"; Value = "511.50" },
    @{ Cell = "E30"; Value = "  +1.23%  " },
    @{ Cell = "D31"; Value = "0.0₃0893" },
    @{ Cell = "E31"; Value = "  -1.06%  " },
    @{ Cell = "D32"; Value = "7.75" },
    @{ Cell = "E32"; Value = "  -0.86%  " },
    @{ Cell = "E33"; Value = "  -2.02%  " },
    @{ Cell = "D34"; Value = "1.77" },
    @{ Cell = "E34"; Value = "  +0.08%  " },
    @{ Cell = "E35"; Value = "  +0.02%  " },
    @{ Cell = "D36"; Value = "161.83" },
    @{ Cell = "E36"; Value = "  -0.43%  " },
    @{ Cell = "E37"; Value = "  -2.62%  " },
    @{ Cell = "D38"; Value = "18.80" },
    @{ Cell = "E38"; Value = "  +2.62%  " },
    @{ Cell = "D39"; Value = "18.69" },
    @{ Cell = "E39"; Value = "  +0.09%  " },
    @{ Cell = "E40"; Value = "  +0.09%  " },
    @{ Cell = "E41"; Value = "  -2.00%  " },
    @{ Cell = "E42"; Value = "  -0.66%  " },
    @{ Cell = "D43"; Value = "4.80" },
    @{ Cell = "E43"; Value = "  -1.04%  " },
    @{ Cell = "E44"; Value = "  -2.54%  " },
    @{ Cell = "B45"; Value = "OKB" },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb" },
    @{ Cell = "D45"; Value = "38.85" },
    @{ Cell = "E45"; Value = "  -0.51%  " },
    @{ Cell = "B46"; Value = "dogwifhat" },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif" },
    @{ Cell = "D46"; Value = "2.33" },
    @{ Cell = "E46"; Value = "  -3.42%  " },
    @{ Cell = "D47"; Value = "149.88" },
    @{ Cell = "E47"; Value = "  +3.58%  " },
    @{ Cell = "E48"; Value = "  +1.73%  " },
    @{ Cell = "D49"; Value = "0.513" },
    @{ Cell = "E49"; Value = "  -0.25%  " },
    @{ Cell = "B50"; Value = "BabyDogeCoin" },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" },
    @{ Cell = "D50"; Value = "0.0₆0251" },
    @{ Cell = "E50"; Value = "  -1.22%  " },
    @{ Cell = "B51"; Value = "Cronos" },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" },
    @{ Cell = "D51"; Value = "0.0738" },
    @{ Cell = "E51"; Value = "  -0.69%  " }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.ClearFormats()
}
